$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd([char]13)
}

# --- 1) "Igår så började vi programmera ..." paragraph: mark "griden" / "gridden"
#        as spell-check proofing errors (w:proofErr spellStart/spellEnd) by
#        splitting the single run into five runs around them. ---
$target1 = "Igår så började vi programmera i java Dennis började med att skapa fönstret och den generella griden, och Emma började med registreringen. Idag ska Emma se till så att formuläret fungerar. Dennis ska fina till den generella gridden och möjligt vi börja med att skapa händelser."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ((Get-ParaText $p) -eq $target1) {
        $xml1 = '<w:p ' + $wNs + '>' +
            '<w:r><w:t xml:space="preserve">Igår så började vi programmera i java Dennis började med att skapa fönstret och den generella </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>griden</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve">, och Emma började med registreringen. Idag ska Emma se till så att formuläret fungerar. Dennis ska fina till den generella </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>gridden</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> och möjligt vi börja med att skapa händelser.</w:t></w:r>' +
            '</w:p>'
        $p.Range.InsertXML($xml1) | Out-Null
        break
    }
}

# --- 2) "Igår så fortsatte Dennis med händelsehanteraren ..." paragraph:
#        mark "Ponuts" / "navbaren" as spell-check proofing errors the same way. ---
$target2 = "Igår så fortsatte Dennis med händelsehanteraren, Emma började med layouten för vecka, månad och dag. Ponuts höll på med navbaren."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ((Get-ParaText $p) -eq $target2) {
        $xml2 = '<w:p ' + $wNs + '>' +
            '<w:r><w:t xml:space="preserve">Igår så fortsatte Dennis med händelsehanteraren, Emma började med layouten för vecka, månad och dag. </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>Ponuts</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> höll på med </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>navbaren</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t>.</w:t></w:r>' +
            '</w:p>'
        $p.Range.InsertXML($xml2) | Out-Null
        break
    }
}

# --- 3) Remove the "_GoBack" bookmark from the last existing paragraph
#        ("Dennis ska göra klart ...") -- it will be re-added, further down,
#        in the new final paragraph. ---
$target3 = "Dennis ska göra klart händelsehanteraren idag, Pontus sa se till att nav baren funkar ordentligt. Emma ska göra klart vyerna och städa i koden."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ((Get-ParaText $p) -eq $target3) {
        $xml3 = '<w:p ' + $wNs + '>' +
            '<w:r><w:t>Dennis ska göra klart händelsehanteraren idag, Pontus sa se till att nav baren funkar ordentligt. Emma ska göra klart vyerna och städa i koden.</w:t></w:r>' +
            '</w:p>'
        $p.Range.InsertXML($xml3) | Out-Null
        break
    }
}

# --- 4) Append the new "2016-12-02" morning-meeting notes at the end of the
#        document, with the "_GoBack" bookmark now living in the new last
#        paragraph. ---
$endRange = $d.Content
$endRange.Collapse(0)
$appendXml = '<w:p ' + $wNs + '><w:r><w:t>2016-12-02</w:t></w:r></w:p>' +
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t xml:space="preserve">Igår: </w:t></w:r>' +
        '<w:r><w:t>Emma gjorde om panelerna, Dennis jobbade med händelser och Pontus</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> jobbade med navigeringen av Emmas paneler.</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t>Idag:</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> Dennis ska fixa layout för händelser så att han blir klar. De andra två ska göra klart det dom gjorde igår.</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
$endRange.InsertXML($appendXml) | Out-Null
